# Update the "想去人数" (column F) counters on the "展览" and "全部类型"
# sheets to reflect the refreshed data output (gh-pages regeneration).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1448
$ws1.Range("F4").Value  = 74
$ws1.Range("F5").Value  = 12085
$ws1.Range("F6").Value  = 4469
$ws1.Range("F11").Value = 2587
$ws1.Range("F14").Value = 60
$ws1.Range("F15").Value = 5235
$ws1.Range("F17").Value = 204
$ws1.Range("F18").Value = 544
$ws1.Range("F19").Value = 11418
$ws1.Range("F20").Value = 11473

# --- Sheet "全部类型" -----------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1448
$ws4.Range("F4").Value  = 74
$ws4.Range("F5").Value  = 12085
$ws4.Range("F6").Value  = 4469
$ws4.Range("F11").Value = 2587
$ws4.Range("F15").Value = 60
$ws4.Range("F16").Value = 5235
$ws4.Range("F18").Value = 204
$ws4.Range("F19").Value = 544
$ws4.Range("F20").Value = 11418
$ws4.Range("F21").Value = 11473
